# Updated cryptos list on Sat Aug  3 22:54:44 UTC 2024 with GitHub Actions
# Refreshes Price (column D) and Volume(1h) (column E) for many coins,
# and reorders Hedera/EnergySwap (rows 44-45) to reflect the new ranking.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @("D2", "60.553.39"),
    @("E2", "  -1.41%  "),
    @("D3", "2.903.59"),
    @("E3", "  -2.50%  "),
    @("E4", "  -0.03%  "),
    @("D5", "529.01"),
    @("E5", "  -2.45%  "),
    @("D6", "143.06"),
    @("E6", "  -6.51%  "),
    @("E7", "  -0.20%  "),
    @("D8", "0.555"),
    @("E8", "  -1.09%  "),
    @("D9", "2.910.66"),
    @("E9", "  -2.30%  "),
    @("E10", "  -2.59%  "),
    @("E11", "  -3.99%  "),
    @("D12", "0.361"),
    @("E12", "  -0.75%  "),
    @("D13", "3.412.86"),
    @("E13", "  -2.55%  "),
    @("E14", "  +1.45%  "),
    @("D15", "60.542.65"),
    @("E15", "  -1.72%  "),
    @("D16", "22.64"),
    @("E16", "  -4.00%  "),
    @("D17", "2.908.70"),
    @("E17", "  -2.61%  "),
    @("E18", "  -3.01%  "),
    @("D19", "5.04"),
    @("E19", "  -0.89%  "),
    @("D20", "11.72"),
    @("E20", "  -1.23%  "),
    @("D21", "364.43"),
    @("E21", "  -6.18%  "),
    @("D22", "6.59"),
    @("E22", "  -0.16%  "),
    @("E23", "  -0.03%  "),
    @("D24", "64.34"),
    @("E24", "  -0.65%  "),
    @("D25", "3.025.23"),
    @("E25", "  -3.07%  "),
    @("E26", "  -2.88%  "),
    @("D27", "0.178"),
    @("E27", "  -3.91%  "),
    @("E28", "  -0.07%  "),
    @("D29", "7.82"),
    @("E29", "  -6.62%  "),
    @("E30", "  -7.91%  "),
    @("E31", "  -0.04%  "),
    @("E32", "  -1.88%  "),
    @("E33", "  -2.98%  "),
    @("D34", "147.43"),
    @("E34", "  -7.50%  "),
    @("D35", "4.36"),
    @("E35", "  -5.54%  "),
    @("D36", "5.59"),
    @("E36", "  -6.74%  "),
    @("E37", "  -5.50%  "),
    @("D38", "1.21"),
    @("E38", "  -5.96%  "),
    @("E39", "  +2.28%  "),
    @("E40", "  -4.70%  "),
    @("D41", "2.330.83"),
    @("E41", "  -4.18%  "),
    @("E42", "  -5.30%  "),
    @("E43", "  -1.82%  "),
    @("E46", "  -0.13%  "),
    @("D47", "4.97"),
    @("E47", "  +1.65%  "),
    @("E48", "  -3.95%  "),
    @("D49", "0.0935"),
    @("E50", "  -1.25%  "),
    @("E51", "  -6.05%  "),
    @("B44", "Hedera"),
    @("C44", "https://coinranking.com/coin/jad286TjB+hedera-hbar"),
    @("D44", "0.0577"),
    @("E44", "  -2.63%  "),
    @("B45", "EnergySwap"),
    @("C45", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"),
    @("D45", "20.65"),
    @("E45", "  -7.03%  ")
)

foreach ($u in $updates) {
    $addr = $u[0]
    $val = $u[1]
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}
